# Yigit Alparslan resume update
# ------------------------------
# The bulk of the source diff is Microsoft Word's own re-save "noise"
# that shows up after the spelling/grammar checker runs over the
# document: <w:proofErr .../> markers get inserted and the runs that
# straddle a checked word/sentence get split into more <w:r> elements,
# but the concatenated visible text of every single one of those spots
# is byte-for-byte identical before and after (verified while analysing
# the diff). There is exactly one place where real, human-authored
# content changed: two new bullet points were appended to the
# "EDUCATION" list, right after the
#   "Bachelor of Science in Computer Engineering ..." line.
# That is the edit this script reproduces.

$d = $word.ActiveDocument

# Find the last Education bullet ("Bachelor of Science in Computer
# Engineering ...") so the insertion point does not depend on hard-coded
# paragraph numbers.
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "Bachelor of Science in Computer Engineering",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Bachelor of Science in Computer Engineering' bullet"
}

$targetIndex = $searchRange.Paragraphs(1).Index

# First new bullet: inherits the ListParagraph / numId=23 bullet
# formatting from the paragraph it is inserted after.
$d.Paragraphs($targetIndex).Range.InsertParagraphAfter()
$d.Paragraphs($targetIndex + 1).Range.Text =
    "Drexel Global Scholar, BS/MS Graduate, Drexel AI Founder, Honors Degree "

# Second new bullet, right after the first one.
$d.Paragraphs($targetIndex + 1).Range.InsertParagraphAfter()
$d.Paragraphs($targetIndex + 2).Range.Text =
    "2015 National Physics Olympiad Silver Medal Winner, Turkey"

Write-Host "Inserted education bullets at paragraphs" ($targetIndex + 1) "and" ($targetIndex + 2)
